$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph ending in "...work-" (it is immediately
# followed by the paragraph that begins with "and by").
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$p1Index = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -eq "I make things - things that work-`r") {
        $p1Index = $i
        break
    }
}
if ($p1Index -eq -1) {
    throw "Could not locate the 'I make things...work-' paragraph"
}
$p2Index = $p1Index + 1

$p2 = $paras.Item($p2Index)
$p2Text = $p2.Range.Text
if (-not ($p2Text.StartsWith("and by"))) {
    throw "Unexpected text in following paragraph: $p2Text"
}

# ------------------------------------------------------------------
# Step 1: extend the first paragraph's (single) run with " and by".
# Because that paragraph currently holds exactly one run, this
# cannot cascade into neighbouring runs.
# ------------------------------------------------------------------
$p1 = $paras.Item($p1Index)
$p1End = $p1.Range.End            # includes the paragraph mark
$insertPoint = $d.Range($p1End - 1, $p1End - 1)
$insertPoint.InsertAfter(" and by")

# ------------------------------------------------------------------
# Step 2: delete the paragraph mark, merging the two paragraphs.
# This operation never coalesces sibling runs.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item($p1Index)
$p1End = $p1.Range.End
$pMark = $d.Range($p1End - 1, $p1End)
$pMark.Delete()

# ------------------------------------------------------------------
# Step 3: the old "and by" run (now the 2nd run of the merged
# paragraph) is redundant and must be removed. Deleting text out of
# a run cascades forward, coalescing every following same-formatted
# run into one - which would destroy the run boundaries the diff
# keeps intact. To stop the cascade we briefly toggle Bold on
# alternating runs that must stay separated, delete the redundant
# text, then restore the formatting.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item($p1Index)
$mergedText = "I make things - things that work- and by"
$dupStart = $p1.Range.Start + $mergedText.Length
$dupEnd = $dupStart + ("and by").Length

$dupCheck = $d.Range($dupStart, $dupEnd)
if ($dupCheck.Text -ne "and by") {
    throw "Unexpected text at dup location: $($dupCheck.Text)"
}

# Runs following the duplicate "and by" run, in order, with lengths:
#   " "(1) "'"(1) "work"(4) "'"(1) " "(1) "I mean work superlatively."(27)
$afterLens = @(1, 1, 4, 1, 1, 27)

$pos = $dupEnd
$ranges = @()
foreach ($len in $afterLens) {
    $ranges += ,@($pos, $pos + $len)
    $pos += $len
}

# Wall off every other run (indices 0,2,4 -> " ", "work", " ") with a
# temporary Bold toggle so they cannot merge with neighbours that
# keep their original formatting.
for ($i = 0; $i -lt $ranges.Length; $i += 2) {
    $b = $ranges[$i]
    $d.Range($b[0], $b[1]).Font.Bold = 1
}

# Remove the duplicate "and by" text.
$dup = $d.Range($dupStart, $dupEnd)
$dup.Delete()

# Restore original formatting on the walled-off runs (positions
# shifted back by the 6 deleted characters).
for ($i = 0; $i -lt $ranges.Length; $i += 2) {
    $b = $ranges[$i]
    $d.Range($b[0] - 6, $b[1] - 6).Font.Bold = 0
}
